$d = $word.ActiveDocument

$pairs = @(
    @("39×35=", "95×14="),
    @("72×69=", "63×96="),
    @("81×36=", "63×78="),
    @("70×60=", "49×58="),
    @("57×27=", "90×78="),
    @("66×27=", "43×28="),
    @("65×93=", "47×86="),
    @("23×66=", "80×74="),
    @("59×26=", "79×42="),
    @("53×62=", "33×31="),
    @("50×93=", "30×77="),
    @("50×36=", "65×31="),
    @("58×77=", "15×16="),
    @("25×22=", "13×39="),
    @("11×27=", "12×33="),
    @("51×63=", "61×53="),
    @("47×56=", "88×20="),
    @("31×50=", "65×36="),
    @("19×65=", "75×71="),
    @("49×31=", "92×25="),
    @("99×23=", "24×32="),
    @("38×88=", "37×96="),
    @("16×59=", "74×50="),
    @("42×67=", "89×79="),
    @("56×12=", "44×17=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
